# Refresh the cached "datetimeFigureOut" Date placeholder text (the
# auto date/time shown in the Slide Master, every Slide Layout, and the
# Notes Master) to 10/8/18 -- mirrors PowerPoint silently re-stamping
# the header/footer date field cache on save.

$p = $ppt.ActivePresentation
$newDate = "10/8/18"
$ppPlaceholderDate = 16

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            $isDatePH = $false
            if ($shp.Type -eq 14) {
                if ($shp.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                    $isDatePH = $true
                }
            }
            if ($isDatePH) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# Slide Master
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every Slide Layout that hangs off the master
for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    Update-DatePlaceholder $layout.Shapes
}

# Notes Master
$notesMaster = $p.NotesMaster
Update-DatePlaceholder $notesMaster.Shapes
